$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60, shifting existing rows 60..142 down to 61..143.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new weekly record.
$ws.Cells.Item(60, 1).Value = 4
$ws.Cells.Item(60, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(60, 3).Value = "Los Lagos"
$ws.Cells.Item(60, 4).Value = 44477
$ws.Cells.Item(60, 5).Value = 10
$ws.Cells.Item(60, 6).Value = "Fruta"
$ws.Cells.Item(60, 7).Value = 100108
$ws.Cells.Item(60, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(60, 9).Value = 100108005
$ws.Cells.Item(60, 10).Value = "Piña"
$ws.Cells.Item(60, 11).Value = "Caramelo"
$ws.Cells.Item(60, 12).Value = "Segunda"
$ws.Cells.Item(60, 13).Value = 100
$ws.Cells.Item(60, 14).Value = 21000
$ws.Cells.Item(60, 15).Value = 21000
$ws.Cells.Item(60, 16).Value = 21000
$ws.Cells.Item(60, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(60, 18).Value = "Ecuador"
$ws.Cells.Item(60, 19).Value = 1500
$ws.Cells.Item(60, 20).Value = 14

# Ensure the date cell keeps the same custom date format used by the rest of column D.
$ws.Cells.Item(60, 4).NumberFormat = $ws.Cells.Item(61, 4).NumberFormat
